# Turkey Super Lig workbook update
# - Swap the B:AC contents between several pairs of adjacent rows (the
#   underlying match records were re-ordered / re-paired in the source data).
# - Replace the contents of rows 319/320 with what used to be rows 323/324
#   (two "future" fixtures that moved up), then drop the now-duplicated
#   trailing rows 321-324 entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($sheet, $rowA, $rowB, $firstCol, $lastCol)

    $rangeA = $sheet.Range("$firstCol$rowA`:$lastCol$rowA")
    $rangeB = $sheet.Range("$firstCol$rowB`:$lastCol$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

# Row pairs whose B:AC data (everything except the id in column A) swap places.
$pairs = @(
    @(171,172),
    @(181,182),
    @(186,187),
    @(197,198),
    @(233,234),
    @(241,242),
    @(250,251),
    @(263,264),
    @(273,274),
    @(281,282)
)

foreach ($pair in $pairs) {
    Swap-RowRange $ws $pair[0] $pair[1] "B" "AC"
}

# Rows 319/320 get overwritten with the data currently sitting in rows 323/324
# (columns B through AA only - these "future fixture" rows have no result
# columns H/I/J and no AB/AC post-match analytics yet).
$ws.Range("B319:AA319").Value2 = $ws.Range("B323:AA323").Value2
$ws.Range("B320:AA320").Value2 = $ws.Range("B324:AA324").Value2

# Rows 321-324 are now redundant (319/320 absorbed 323/324's data, and
# 321/322 are dropped outright) - remove them so the sheet ends at row 320.
$ws.Rows("321:324").Delete()
